# Trade #103 (row 104 / "All Trades", row 24 / "momentum") closes as an
# early_exit, and a brand-new trade #132 ("MarketMaking") opens, appended as
# a new last row on both "All Trades" and "MarketMaking". Summary and
# Strategy Status roll-up figures move accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.36   # Current Capital
$summary.Range("B4").Value = 0.47      # Total P&L $
$summary.Range("B6").Value = 103       # Total Trades
$summary.Range("B7").Value = 49        # Winning Trades
$summary.Range("B9").Value = 47.57     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet, "momentum" row (row 11)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.29000000000001
$status.Range("D11").Value = 23
$status.Range("E11").Value = -0.71
$status.Range("F11").Value = -0.71
$status.Range("G11").Value = 26.09

# ---------------------------------------------------------------------
# All Trades sheet, row 104 -> trade #103 ("momentum") closes early
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G104").Value = 0.85
$allTrades.Range("H104").Value = "CLOSED"
$allTrades.Range("I104").Value = 3.6585
$allTrades.Range("J104").Value = 0.03
$allTrades.Range("K104").Value = 99.29000000000001
$allTrades.Range("L104").Value = "early_exit"
$allTrades.Range("M104").Value = 0.13

# ---------------------------------------------------------------------
# momentum sheet, row 24 -> same trade #103, mirrored
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("G24").Value = 0.85
$momentum.Range("H24").Value = "CLOSED"
$momentum.Range("I24").Value = 3.6585
$momentum.Range("J24").Value = 0.03
$momentum.Range("K24").Value = 99.29000000000001
$momentum.Range("P24").Value = "early_exit"
$momentum.Range("Q24").Value = 0.13

# ---------------------------------------------------------------------
# New trade #132 ("MarketMaking", still OPEN) appended to "All Trades"
# as row 133, and to the "MarketMaking" strategy sheet as row 53.
# Date/Time text columns get an explicit Text format first so the COM
# layer doesn't silently coerce the "2026-02-18" looking string into a
# date serial number.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- All Trades!A133:Q133 ---
$allTrades.Range("A133").Value = 132
Set-TextValue $allTrades.Range("B133") "2026-02-18"
Set-TextValue $allTrades.Range("C133") "00:27:47"
$allTrades.Range("D133").Value = "MarketMaking"
$allTrades.Range("E133").Value = "DOWN"
$allTrades.Range("F133").Value = 0.82
# G133 (Exit Price) intentionally left blank - trade is still OPEN
$allTrades.Range("H133").Value = "OPEN"
$allTrades.Range("I133").Value = 0
$allTrades.Range("J133").Value = 0
$allTrades.Range("K133").Value = 99.47967800952271
# L133 (Exit Reason) intentionally left blank - trade is still OPEN
$allTrades.Range("M133").Value = 0
$allTrades.Range("N133").Value = 0
$allTrades.Range("O133").Value = 0
$allTrades.Range("P133").Value = 0.65
$allTrades.Range("Q133").Value = "Wide spread capture: 392 bps vs avg 299 bps"

# --- MarketMaking!A53:Q53 (same trade, strategy-sheet column layout) ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A53").Value = 132
Set-TextValue $marketMaking.Range("B53") "2026-02-18"
Set-TextValue $marketMaking.Range("C53") "00:27:47"
$marketMaking.Range("D53").Value = "MarketMaking"
$marketMaking.Range("E53").Value = "DOWN"
$marketMaking.Range("F53").Value = 0.82
# G53 (Exit Price) intentionally left blank - trade is still OPEN
$marketMaking.Range("H53").Value = "OPEN"
$marketMaking.Range("I53").Value = 0
$marketMaking.Range("J53").Value = 0
$marketMaking.Range("K53").Value = 99.47967800952271
$marketMaking.Range("L53").Value = 0
$marketMaking.Range("M53").Value = 0
$marketMaking.Range("N53").Value = 0.65
$marketMaking.Range("O53").Value = "Wide spread capture: 392 bps vs avg 299 bps"
# P53 (Exit Reason) intentionally left blank - trade is still OPEN
$marketMaking.Range("Q53").Value = 0

Write-Host "Edit applied."
